$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add data to the stable upright experiment statistics
$ws.Range("H7").Value = 2
$ws.Range("H8").Value = 5

# Update the active cell selection as recorded in the workbook view state
$ws.Range("H14").Select()
